$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.683.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "

# Row 3
$ws.Range("D3").Value = "'1.849.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").Value = "'312.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("D7").Value = "'0.4281"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "

# Row 8
$ws.Range("E8").Value = "  -1.51%  "

# Row 9
$ws.Range("D9").Value = "'0.07296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "

# Row 10
$ws.Range("D10").Value = "'0.8714"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "

# Row 11
$ws.Range("D11").Value = "'20.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "

# Row 12
$ws.Range("D12").Value = "'1.907.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.40%  "

# Row 13
$ws.Range("D13").Value = "'6.555"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.14%  "

# Row 14
$ws.Range("D14").Value = "'5.331"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "

# Row 15
$ws.Range("D15").Value = "'0.07019"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "

# Row 16
$ws.Range("D16").Value = "'1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "

# Row 17
$ws.Range("D17").Value = "'79.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("D18").Value = "'0.000008948"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("D19").Value = "'1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "

# Row 20
$ws.Range("D20").Value = "'15.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.92%  "

# Row 21
$ws.Range("D21").Value = "'27.697.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("D22").Value = "'5.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "

# Row 23
$ws.Range("D23").Value = "'10.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.90%  "

# Row 24
$ws.Range("D24").Value = "'2.074.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

# Row 25
$ws.Range("E25").Value = "  +2.46%  "

# Row 26
$ws.Range("D26").Value = "'155.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.13%  "

# Row 27
$ws.Range("D27").Value = "'18.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.25%  "

# Row 28
$ws.Range("D28").Value = "'120.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "

# Row 29
$ws.Range("D29").Value = "'5.279"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "

# Row 30
$ws.Range("D30").Value = "'1.875"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.93%  "

# Row 31
$ws.Range("E31").Value = "  -0.10%  "

# Row 32
$ws.Range("D32").Value = "'0.7577"
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'2.970"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.03%  "

# Row 34
$ws.Range("D34").Value = "'4.513"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.40%  "

# Row 35
$ws.Range("E35").Value = "  +2.49%  "

# Row 36
$ws.Range("E36").Value = "  +0.55%  "

# Row 37
$ws.Range("D37").Value = "'0.05426"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.00%  "

# Row 38
$ws.Range("D38").Value = "'1.100"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "

# Row 39
$ws.Range("E39").Value = "  -0.29%  "

# Row 40
$ws.Range("D40").Value = "'2.827"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "

# Row 41
$ws.Range("D41").Value = "'0.1666"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.77%  "

# Row 42
$ws.Range("D42").Value = "'0.5073"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "

# Row 43
$ws.Range("D43").Value = "'6.617"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.88%  "

# Row 44
$ws.Range("D44").Value = "'8.414"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'106.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.14%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.06539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "

# Row 47
$ws.Range("E47").Value = "  +0.29%  "

# Row 48
$ws.Range("D48").Value = "'0.4652"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "

# Row 49
$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "

# Row 50
$ws.Range("D50").Value = "'1.628"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.51%  "

# Row 51
$ws.Range("D51").Value = "'1.785"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
